$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Datos")
$ws2 = $wb.Worksheets.Item("Hojita 2")

# New header columns
$ws1.Range("D1").Value = "Edad"
$ws1.Range("E1").Value = "Estado"
$ws1.Range("F1").Value = "Altura"

# Row 2 (Camilo Arguello) - keep name/surname, add new data
$ws1.Range("D2").Value = 24
$ws1.Range("E2").Value = "Soltero"
$ws1.Range("F2").Value = 1.7

# Row 3 - change name/surname to Juli Ruiz, add new data
$ws1.Range("B3").Value = "Juli"
$ws1.Range("C3").Value = "Ruiz"
$ws1.Range("D3").Value = 23
$ws1.Range("E3").Value = "Soltera"
$ws1.Range("F3").Value = 1.71

# Row 4 - change name/surname to Andres Rincon, add new data
$ws1.Range("B4").Value = "Andres"
$ws1.Range("C4").Value = "Rincon"
$ws1.Range("D4").Value = 22
$ws1.Range("E4").Value = "Casado"
$ws1.Range("F4").Value = 1.5

# Update the selections / active sheet to match the final view state
$ws2.Range("B2").Select()

$ws1.Activate()
$ws1.Range("G4").Select()
